$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Add a new entry on row 20: date, start time, end time, and category.
# The duration column (D) already contains a filled-down shared formula
# that will pick up the new row automatically.
$ws.Range("A20").Value = 45488
$ws.Range("B20").Value = 0.66666666666666663
$ws.Range("C20").Value = 0.95277777777777783
$ws.Range("E20").Value = "Analyse et état de l'art"

# Move the selection to where the user ended up after the edit.
$ws.Range("H17").Select()
